# Auto-generated: update TPM-derived NATMI metrics for F2-Itga2b sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.754521
$ws.Range("H2").Value = 2.263563
$ws.Range("I2").Value = 0.2768403531129761
$ws.Range("J2").Value = 0.2768403531129761
$ws.Range("M2").Value = 1.646992
$ws.Range("N2").Value = 4.940976
$ws.Range("O2").Value = 0.2071783517404009
$ws.Range("P2").Value = 0.2071783517404009
$ws.Range("Q2").Value = 1.242690050832
$ws.Range("R2").Value = 11.184210457488
$ws.Range("S2").Value = 0.05735532805317697
$ws.Range("T2").Value = 0.05735532805317697
# Row 3
$ws.Range("G3").Value = 0.754521
$ws.Range("H3").Value = 2.263563
$ws.Range("I3").Value = 0.2768403531129761
$ws.Range("J3").Value = 0.2768403531129761
$ws.Range("O3").Value = 0.4685125322965616
$ws.Range("P3").Value = 0.4685125322965616
$ws.Range("Q3").Value = 2.810215728063
$ws.Range("R3").Value = 25.291941552567
$ws.Range("S3").Value = 0.1297031748788347
$ws.Range("T3").Value = 0.1297031748788347
# Row 4
$ws.Range("G4").Value = 0.754521
$ws.Range("H4").Value = 2.263563
$ws.Range("I4").Value = 0.2768403531129761
$ws.Range("J4").Value = 0.2768403531129761
$ws.Range("M4").Value = 2.284352333333333
$ws.Range("N4").Value = 6.853057
$ws.Range("O4").Value = 0.2873531572796583
$ws.Range("P4").Value = 0.2873531572796583
$ws.Range("Q4").Value = 1.723591806899
$ws.Range("R4").Value = 15.512326262091
$ws.Range("S4").Value = 0.07955094952942916
$ws.Range("T4").Value = 0.07955094952942916
# Row 5
$ws.Range("G5").Value = 0.754521
$ws.Range("H5").Value = 2.263563
$ws.Range("I5").Value = 0.2768403531129761
$ws.Range("J5").Value = 0.2768403531129761
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2937863333333333
$ws.Range("N5").Value = 0.881359
$ws.Range("O5").Value = 0.03695595868337916
$ws.Range("P5").Value = 0.03695595868337916
$ws.Range("Q5").Value = 0.221667958013
$ws.Range("R5").Value = 1.995011622117
$ws.Range("S5").Value = 0.01023090065153524
$ws.Range("T5").Value = 0.01023090065153524
# Row 6
$ws.Range("G6").Value = 0.9731926666666667
$ws.Range("I6").Value = 0.3570728998754956
$ws.Range("J6").Value = 0.3570728998754956
$ws.Range("M6").Value = 1.646992
$ws.Range("N6").Value = 4.940976
$ws.Range("O6").Value = 0.2071783517404009
$ws.Range("P6").Value = 0.2071783517404009
$ws.Range("Q6").Value = 1.602840536458667
$ws.Range("R6").Value = 14.425564828128
$ws.Range("S6").Value = 0.07397777484737041
$ws.Range("T6").Value = 0.07397777484737041
# Row 7
$ws.Range("G7").Value = 0.9731926666666667
$ws.Range("I7").Value = 0.3570728998754956
$ws.Range("J7").Value = 0.3570728998754956
$ws.Range("O7").Value = 0.4685125322965616
$ws.Range("P7").Value = 0.4685125322965616
$ws.Range("S7").Value = 0.1672931285351451
$ws.Range("T7").Value = 0.1672931285351451
# Row 8
$ws.Range("G8").Value = 0.9731926666666667
$ws.Range("I8").Value = 0.3570728998754956
$ws.Range("J8").Value = 0.3570728998754956
$ws.Range("M8").Value = 2.284352333333333
$ws.Range("N8").Value = 6.853057
$ws.Range("O8").Value = 0.2873531572796583
$ws.Range("P8").Value = 0.2873531572796583
$ws.Range("Q8").Value = 2.223114938882889
$ws.Range("R8").Value = 20.008034449946
$ws.Range("S8").Value = 0.102606025158227
$ws.Range("T8").Value = 0.102606025158227
# Row 9
$ws.Range("G9").Value = 0.9731926666666667
$ws.Range("I9").Value = 0.3570728998754956
$ws.Range("J9").Value = 0.3570728998754956
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2937863333333333
$ws.Range("N9").Value = 0.881359
$ws.Range("O9").Value = 0.03695595868337916
$ws.Range("P9").Value = 0.03695595868337916
$ws.Range("Q9").Value = 0.2859107051668889
$ws.Range("R9").Value = 2.573196346502
$ws.Range("S9").Value = 0.0131959713347532
$ws.Range("T9").Value = 0.0131959713347532
# Row 10
$ws.Range("G10").Value = 0.7824410000000001
$ws.Range("H10").Value = 2.347323
$ws.Range("I10").Value = 0.2870844452706686
$ws.Range("J10").Value = 0.2870844452706686
$ws.Range("M10").Value = 1.646992
$ws.Range("N10").Value = 4.940976
$ws.Range("O10").Value = 0.2071783517404009
$ws.Range("P10").Value = 0.2071783517404009
$ws.Range("Q10").Value = 1.288674067472
$ws.Range("R10").Value = 11.598066607248
$ws.Range("S10").Value = 0.05947768218148446
$ws.Range("T10").Value = 0.05947768218148447
# Row 11
$ws.Range("G11").Value = 0.7824410000000001
$ws.Range("H11").Value = 2.347323
$ws.Range("I11").Value = 0.2870844452706686
$ws.Range("J11").Value = 0.2870844452706686
$ws.Range("O11").Value = 0.4685125322965616
$ws.Range("P11").Value = 0.4685125322965616
$ws.Range("Q11").Value = 2.914203851823001
$ws.Range("R11").Value = 26.22783466640701
$ws.Range("S11").Value = 0.1345026604367146
$ws.Range("T11").Value = 0.1345026604367146
# Row 12
$ws.Range("G12").Value = 0.7824410000000001
$ws.Range("H12").Value = 2.347323
$ws.Range("I12").Value = 0.2870844452706686
$ws.Range("J12").Value = 0.2870844452706686
$ws.Range("M12").Value = 2.284352333333333
$ws.Range("N12").Value = 6.853057
$ws.Range("O12").Value = 0.2873531572796583
$ws.Range("P12").Value = 0.2873531572796583
$ws.Range("Q12").Value = 1.787370924045667
$ws.Range("R12").Value = 16.086338316411
$ws.Range("S12").Value = 0.08249462175440588
$ws.Range("T12").Value = 0.08249462175440589
# Row 13
$ws.Range("G13").Value = 0.7824410000000001
$ws.Range("H13").Value = 2.347323
$ws.Range("I13").Value = 0.2870844452706686
$ws.Range("J13").Value = 0.2870844452706686
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.2937863333333333
$ws.Range("N13").Value = 0.881359
$ws.Range("O13").Value = 0.03695595868337916
$ws.Range("P13").Value = 0.03695595868337916
$ws.Range("Q13").Value = 0.2298704724396667
$ws.Range("R13").Value = 2.068834251957
$ws.Range("S13").Value = 0.01060948089806365
$ws.Range("T13").Value = 0.01060948089806366
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2153186666666667
$ws.Range("H14").Value = 0.645956
$ws.Range("I14").Value = 0.07900230174085969
$ws.Range("J14").Value = 0.07900230174085969
$ws.Range("M14").Value = 1.646992
$ws.Range("N14").Value = 4.940976
$ws.Range("O14").Value = 0.2071783517404009
$ws.Range("P14").Value = 0.2071783517404009
$ws.Range("Q14").Value = 0.3546281214506666
$ws.Range("R14").Value = 3.191653093056
$ws.Range("S14").Value = 0.01636756665836912
$ws.Range("T14").Value = 0.01636756665836912
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2153186666666667
$ws.Range("H15").Value = 0.645956
$ws.Range("I15").Value = 0.07900230174085969
$ws.Range("J15").Value = 0.07900230174085969
$ws.Range("O15").Value = 0.4685125322965616
$ws.Range("P15").Value = 0.4685125322965616
$ws.Range("Q15").Value = 0.801955019956
$ws.Range("R15").Value = 7.217595179604
$ws.Range("S15").Value = 0.03701356844586723
$ws.Range("T15").Value = 0.03701356844586723
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2153186666666667
$ws.Range("H16").Value = 0.645956
$ws.Range("I16").Value = 0.07900230174085969
$ws.Range("J16").Value = 0.07900230174085969
$ws.Range("M16").Value = 2.284352333333333
$ws.Range("N16").Value = 6.853057
$ws.Range("O16").Value = 0.2873531572796583
$ws.Range("P16").Value = 0.2873531572796583
$ws.Range("Q16").Value = 0.4918636986102222
$ws.Range("R16").Value = 4.426773287492
$ws.Range("S16").Value = 0.02270156083759628
$ws.Range("T16").Value = 0.02270156083759628
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.2153186666666667
$ws.Range("H17").Value = 0.645956
$ws.Range("I17").Value = 0.07900230174085969
$ws.Range("J17").Value = 0.07900230174085969
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.2937863333333333
$ws.Range("N17").Value = 0.881359
$ws.Range("O17").Value = 0.03695595868337916
$ws.Range("P17").Value = 0.03695595868337916
$ws.Range("Q17").Value = 0.06325768157822222
$ws.Range("R17").Value = 0.5693191342039999
$ws.Range("S17").Value = 0.002919605799027064
$ws.Range("T17").Value = 0.002919605799027065
